$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tblIngredients")

# Row 74 currently lists "sugar" / "Grocery"; the new layout inserts a
# "pizza" / "Grocery" row after it, so "sugar" becomes "Check" to make
# room for the newly appended ingredient row below it.
$ws.Range("B74").Value = "Check"

# New row 75: pizza / Grocery, with zeroed nutrition columns.
$ws.Range("A75").Value = "pizza"
$ws.Range("B75").Value = "Grocery"
$ws.Range("C75:G75").Value = 0

# Move the active selection to match the new layout.
$ws.Range("B74").Select() | Out-Null
